$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New D2 cell (dataset 1 now has a size/link placeholder in D column)
$ws.Range("D2").Value = '>100M'

# Row 10: (9) CN A-shares price history dataset
$ws.Range("A10").Value = '(9)A股3000只股票历史价格'
$ws.Range("B10").Value = '19910101-20170620'
$ws.Range("C10").Value = '161M'
$ws.Range("D10").Value = '>100M'
$ws.Range("E10").Value = '链接: https://pan.baidu.com/s/1nvPZK6X 密码: i5ae'
$ws.Range("F10").Value = '复权后的价格数据。格式为[date, open, high, close, low, ?, volume]'
$ws.Range("G10").Value = '由白天煜同学提供。使用Python的Tushare包获取。'

# Row 11: (10) Nasdaq 6000 stocks 10y dataset
$ws.Range("A11").Value = '(10)Nasdaq6000只股票10年数据'
$ws.Range("B11").Value = '20070618-20170616'
$ws.Range("C11").Value = '125M'
$ws.Range("D11").Value = '>100M'
$ws.Range("E11").Value = '链接: https://pan.baidu.com/s/1nvhu3CP 密码: sfwx'
$ws.Range("F11").Value = 'Nasdaq市场股票十年的价格数据。每只股票对应一个文件。文件每行格式为：[date + "\t" + open + "\t" + high + "\t" + low + "\t" + close + "\t" + volume]。【注意】部分股票价格不全（那些只有65天数据的），可设置爬虫程序，让程序在请求数据后等待更长时间，以获得完全的数据。'
$ws.Range("G11").Value = '纳斯达克官网http://www.nasdaq.com/symbol/[symbol]/historical'

# Update view: scroll so column F is at the left edge, then select F13
$excel.Goto($ws.Range("F1"), $true)
$ws.Range("F13").Select()